$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM data values for rows 2-21 (residues), columns B-K (positions 0-9)
$data = @(
        @(-16.01497939028436, 0.5121463070642013, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436),
        @(-16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, 4.321907397664166, -16.01497939028436, -16.01497939028436),
        @(-16.01497939028436, 0.5516247132411835, 0.8791870621349428, -16.01497939028436, 3.818486729166845, -16.01497939028436, 1.444268101358496, -16.01497939028436, 2.958227481256104, -16.01497939028436),
        @(-16.01497939028436, 0.568764500127538, -16.01497939028436, -16.01497939028436, -16.01497939028436, 2.842287878237776, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436),
        @(-16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436),
        @(3.156062989131342, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436),
        @(-16.01497939028436, -16.01497939028436, -16.01497939028436, 2.204420773365121, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436),
        @(3.470639651546403, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436),
        @(-16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, 2.131713574930489),
        @(-16.01497939028436, -16.01497939028436, -16.01497939028436, 2.310997954862974, -16.01497939028436, 2.65554201475793, -16.01497939028436, -16.01497939028436, -16.01497939028436, 1.490281683864032),
        @(-16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436),
        @(-16.01497939028436, -16.01497939028436, -16.01497939028436, 1.852039874612524, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, 1.760097136550888, 1.711165578522448),
        @(-16.01497939028436, -16.01497939028436, 1.818831046872734, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, 2.642482121322741),
        @(-16.01497939028436, -16.01497939028436, 0.1880817888926133, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436),
        @(-16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, 1.152840537208617, -16.01497939028436),
        @(-16.01497939028436, 1.204614568680213, 0.8270885352215253, -16.01497939028436, -16.01497939028436, -16.01497939028436, 2.285127348933623, -16.01497939028436, 2.151070879306556, -16.01497939028436),
        @(-16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, -16.01497939028436, 2.277418639193094, -16.01497939028436, 1.121541183627227, -16.01497939028436),
        @(-16.01497939028436, -16.01497939028436, 2.674298609564294, -16.01497939028436, -16.01497939028436, -16.01497939028436, 1.505247934401049, -16.01497939028436, -16.01497939028436, -16.01497939028436),
        @(-16.01497939028436, 2.770731682303062, 2.415722246527841, -16.01497939028436, 2.558582744089955, -16.01497939028436, 1.128533370245382, -16.01497939028436, -16.01497939028436, 1.717936594359474),
        @(-16.01497939028436, 2.699403654533488, -16.01497939028436, 2.769454172393043, -16.01497939028436, 2.706497008944138, 1.339863989144134, -16.01497939028436, -16.01497939028436, -16.01497939028436)
    )

for ($r = 0; $r -lt 20; $r++) {
    for ($c = 0; $c -lt 10; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $data[$r][$c]
    }
}
